$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.875.12'
$ws.Range("E2").Value = '  +2.43%  '

$ws.Range("D3").Value = '3.808.49'
$ws.Range("E3").Value = '  +0.86%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Formula = "'702.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +11.47%  '

$ws.Range("D6").Formula = "'173.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.51%  '

$ws.Range("D7").Value = '3.807.31'
$ws.Range("E7").Value = '  +0.88%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("E9").Value = '  +1.00%  '

$ws.Range("E10").Value = '  +2.37%  '

$ws.Range("E11").Value = '  +13.06%  '

$ws.Range("E12").Value = '  +0.58%  '

$ws.Range("D13").Formula = "'0.0000255"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.53%  '

$ws.Range("D14").Formula = "'36.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.90%  '

$ws.Range("D15").Value = '4.449.06'
$ws.Range("E15").Value = '  +0.91%  '

$ws.Range("D16").Value = '3.817.34'
$ws.Range("E16").Value = '  +1.08%  '

$ws.Range("D17").Value = '70.889.46'
$ws.Range("E17").Value = '  +2.45%  '

$ws.Range("D18").Formula = "'17.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.58%  '

$ws.Range("E19").Value = '  +2.94%  '

$ws.Range("E20").Value = '  +0.32%  '

$ws.Range("D21").Formula = "'11.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +18.84%  '

$ws.Range("D22").Formula = "'480.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.50%  '

$ws.Range("D23").Formula = "'0.716"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.19%  '

$ws.Range("D24").Formula = "'83.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.14%  '

$ws.Range("D25").Formula = "'0.0000146"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.71%  '

$ws.Range("D26").Formula = "'12.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.79%  '

$ws.Range("E27").Value = '  +1.16%  '

$ws.Range("D28").Formula = "'10.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.91%  '

$ws.Range("D29").Value = '3.959.59'
$ws.Range("E29").Value = '  +0.88%  '

$ws.Range("E30").Value = '  -0.05%  '

$ws.Range("D31").Formula = "'3.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +16.03%  '

$ws.Range("E32").Value = '  +1.80%  '

$ws.Range("D33").Formula = "'7.52"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.77%  '

$ws.Range("D34").Formula = "'29.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.61%  '

$ws.Range("D35").Formula = "'0.179"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.84%  '

$ws.Range("E36").Value = '  +2.60%  '

$ws.Range("D37").Formula = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.16%  '

$ws.Range("D38").Value = '3.758.41'
$ws.Range("E38").Value = '  +0.81%  '

$ws.Range("E39").Value = '  +1.55%  '

$ws.Range("D40").Formula = "'3.55"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.14%  '

$ws.Range("D41").Formula = "'5.99"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.27%  '

$ws.Range("E42").Value = '  +24.71%  '

$ws.Range("D43").Formula = "'2.20"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +13.15%  '

$ws.Range("D44").Formula = "'0.967"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.44%  '

$ws.Range("E45").Value = '  +0.00%  '

$ws.Range("E46").Value = '  +0.00%  '

$ws.Range("D47").Formula = "'45.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.73%  '

$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").Formula = "'160.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.32%  '

$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D49").Formula = "'49.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.36%  '

$ws.Range("E50").Value = '  -0.18%  '

$ws.Range("E51").Value = '  +1.48%  '
